$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worked examples")
$ws.Activate()

# "Excess of Loss 2 Occ Lim" (the 2nd reinsurance inuring layer's occurrence limit)
# changed from 100 to 200
$ws.Range("C23").Value = 200

# Leave the selection where the author left it when saving
$ws.Range("C47:N47").Select() | Out-Null
